$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3057
$ws.Range("E2").Value = 30
$ws.Range("F2").Value = 30
$ws.Range("G2").Value = 39
$ws.Range("H2").Value = 42
$ws.Range("I2").Value = 43
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 2327
$ws.Range("L2").Value = 970
$ws.Range("M2").Value = 1357
$ws.Range("N2").Value = 1358
$ws.Range("O2").Value = -1
$ws.Range("P2").Value = 58
$ws.Range("Q2").Value = 300
$ws.Range("R2").Value = -188
$ws.Range("S2").Value = -34
$ws.Range("T2").Value = 252
$ws.Range("U2").Value = 47
$ws.Range("V2").Value = 337
$ws.Range("W2").Value = 0.97
$ws.Range("X2").Value = 1.38
$ws.Range("Y2").Value = 3.17
$ws.Range("Z2").Value = 1.84
$ws.Range("AA2").Value = 71.47
$ws.Range("AB2").Value = 2317.59
$ws.Range("AC2").Value = 371
$ws.Range("AD2").Value = 20.47
$ws.Range("AE2").Value = 12014
$ws.Range("AF2").Value = 0.63
$ws.Range("AG2").Value = 175
$ws.Range("AH2").Value = 2.31
$ws.Range("AI2").Value = 46.37
$ws.Range("AJ2").Value = 11500000

# Row 3
$ws.Range("D3").Value = 2967
$ws.Range("E3").Value = 40
$ws.Range("F3").Value = 40
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 45
$ws.Range("I3").Value = 46
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 2594
$ws.Range("L3").Value = 1210
$ws.Range("M3").Value = 1384
$ws.Range("N3").Value = 1384
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 58
$ws.Range("Q3").Value = 163
$ws.Range("R3").Value = -221
$ws.Range("S3").Value = 56
$ws.Range("T3").Value = 262
$ws.Range("U3").Value = -99
$ws.Range("V3").Value = 412
$ws.Range("W3").Value = 1.34
$ws.Range("X3").Value = 1.52
$ws.Range("Y3").Value = 3.33
$ws.Range("Z3").Value = 1.83
$ws.Range("AA3").Value = 87.47
$ws.Range("AB3").Value = 2350.32
$ws.Range("AC3").Value = 397
$ws.Range("AD3").Value = 19.81
$ws.Range("AE3").Value = 12246
$ws.Range("AF3").Value = 0.64
$ws.Range("AG3").Value = 175
$ws.Range("AH3").Value = 2.22
$ws.Range("AI3").Value = 43.28
$ws.Range("AJ3").Value = 11500000

# Row 4
$ws.Range("D4").Value = 3232
$ws.Range("E4").Value = 125
$ws.Range("F4").Value = 125
$ws.Range("G4").Value = 143
$ws.Range("H4").Value = 114
$ws.Range("I4").Value = 114
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2989
$ws.Range("L4").Value = 1494
$ws.Range("M4").Value = 1494
$ws.Range("N4").Value = 1494
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 58
$ws.Range("Q4").Value = 516
$ws.Range("R4").Value = -347
$ws.Range("S4").Value = 28
$ws.Range("T4").Value = 385
$ws.Range("U4").Value = 131
$ws.Range("V4").Value = 461
$ws.Range("W4").Value = 3.88
$ws.Range("X4").Value = 3.52
$ws.Range("Y4").Value = 7.9
$ws.Range("Z4").Value = 4.07
$ws.Range("AA4").Value = 100.01
$ws.Range("AB4").Value = 2509.12
$ws.Range("AC4").Value = 989
$ws.Range("AD4").Value = 9.71
$ws.Range("AE4").Value = 13223
$ws.Range("AF4").Value = 0.73
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 2.08
$ws.Range("AI4").Value = 19.87
$ws.Range("AJ4").Value = 11500000

# Row 5
$ws.Range("D5").Value = 3397
$ws.Range("E5").Value = 203
$ws.Range("F5").Value = 203
$ws.Range("G5").Value = 239
$ws.Range("H5").Value = 171
$ws.Range("I5").Value = 171
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 2576
$ws.Range("L5").Value = 1014
$ws.Range("M5").Value = 1562
$ws.Range("N5").Value = 1562
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 58
$ws.Range("Q5").Value = 381
$ws.Range("R5").Value = -167
$ws.Range("S5").Value = -173
$ws.Range("T5").Value = 141
$ws.Range("U5").Value = 240
$ws.Range("V5").Value = 310
$ws.Range("W5").Value = 5.99
$ws.Range("X5").Value = 5.04
$ws.Range("Y5").Value = 11.19
$ws.Range("Z5").Value = 6.15
$ws.Range("AA5").Value = 64.9
$ws.Range("AB5").Value = 2759.51
$ws.Range("AC5").Value = 1487
$ws.Range("AD5").Value = 7.87
$ws.Range("AE5").Value = 13826
$ws.Range("AF5").Value = 0.85
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 1.71
$ws.Range("AI5").Value = 13.21
$ws.Range("AJ5").Value = 11500000

# Row 6
$ws.Range("D6").Value = 3607
$ws.Range("E6").Value = 261
$ws.Range("F6").Value = 261
$ws.Range("G6").Value = 277
$ws.Range("H6").Value = 187
$ws.Range("I6").Value = 187
$ws.Range("K6").Value = 2848
$ws.Range("L6").Value = 1100
$ws.Range("M6").Value = 1748
$ws.Range("N6").Value = 1748
$ws.Range("P6").Value = 58
$ws.Range("Q6").Value = 535
$ws.Range("R6").Value = -395
$ws.Range("S6").Value = -97
$ws.Range("T6").Value = 233
$ws.Range("U6").Value = 302
$ws.Range("V6").Value = 235
$ws.Range("W6").Value = 7.23
$ws.Range("X6").Value = 5.19
$ws.Range("Y6").Value = 11.32
$ws.Range("Z6").Value = 6.91
$ws.Range("AA6").Value = 62.94
$ws.Range("AB6").Value = 3033.78
$ws.Range("AC6").Value = 1630
$ws.Range("AD6").Value = 10.52
$ws.Range("AE6").Value = 15466
$ws.Range("AF6").Value = 1.11
$ws.Range("AG6").Value = 200
$ws.Range("AH6").Value = 1.17
$ws.Range("AI6").Value = 12.06
$ws.Range("AJ6").Value = 11500000

# Clear rows 7-9 (keep only A, B, C columns)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
